$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GitHub commit numbers / day counts
$ws.Range("C8").Value = 14
$ws.Range("C9").Value = 22
$ws.Range("C11").Value = 8
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 5
$ws.Range("C26").Value = 10

# Fill in "No" for the Admin options (rows 34-50, column C), previously blank
for ($r = 34; $r -le 50; $r++) {
    $ws.Cells.Item($r, 3).Value = "No"
}

# Update the active selection to match the saved view state
$ws.Range("C19").Select()
